$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Vegfc"
$ws.Range("C2").Value = "Vipr2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 3.263122
$ws.Range("H2").Value = 9.789366
$ws.Range("I2").Value = 0.3531375780718168
$ws.Range("J2").Value = 0.3531375780718168
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.07178566666666666
$ws.Range("N2").Value = 0.215357
$ws.Range("O2").Value = 0.00516242483834057
$ws.Range("P2").Value = 0.00516242483834057
$ws.Range("Q2").Value = 0.2342453881846666
$ws.Range("R2").Value = 2.108208493662
$ws.Range("S2").Value = 0.001823046204389379
$ws.Range("T2").Value = 0.001823046204389379

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Vegfc"
$ws.Range("C3").Value = "Vipr2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 3.263122
$ws.Range("H3").Value = 9.789366
$ws.Range("I3").Value = 0.3531375780718168
$ws.Range("J3").Value = 0.3531375780718168
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 6.804012333333334
$ws.Range("N3").Value = 20.412037
$ws.Range("O3").Value = 0.4893066248597758
$ws.Range("P3").Value = 0.4893066248597758
$ws.Range("Q3").Value = 22.20232233317133
$ws.Range("R3").Value = 199.820900998542
$ws.Range("S3").Value = 0.1727925564374763
$ws.Range("T3").Value = 0.1727925564374763

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Vegfc"
$ws.Range("C4").Value = "Vipr2"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 3.263122
$ws.Range("H4").Value = 9.789366
$ws.Range("I4").Value = 0.3531375780718168
$ws.Range("J4").Value = 0.3531375780718168
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.029003
$ws.Range("N4").Value = 0.087009
$ws.Range("O4").Value = 0.002085734026566003
$ws.Range("P4").Value = 0.002085734026566003
$ws.Range("Q4").Value = 0.09464032736599999
$ws.Range("R4").Value = 0.851762946294
$ws.Range("S4").Value = 0.0007365510626434966
$ws.Range("T4").Value = 0.0007365510626434966

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Vegfc"
$ws.Range("C5").Value = "Vipr2"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 3.263122
$ws.Range("H5").Value = 9.789366
$ws.Range("I5").Value = 0.3531375780718168
$ws.Range("J5").Value = 0.3531375780718168
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 7.000615333333333
$ws.Range("N5").Value = 21.001846
$ws.Range("O5").Value = 0.5034452162753175
$ws.Range("P5").Value = 0.5034452162753175
$ws.Range("Q5").Value = 22.84386190773733
$ws.Range("R5").Value = 205.594757169636
$ws.Range("S5").Value = 0.1777854243673076
$ws.Range("T5").Value = 0.1777854243673076

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Vegfc"
$ws.Range("C6").Value = "Vipr2"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 4.367310666666667
$ws.Range("H6").Value = 13.101932
$ws.Range("I6").Value = 0.4726337266929886
$ws.Range("J6").Value = 0.4726337266929886
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.07178566666666666
$ws.Range("N6").Value = 0.215357
$ws.Range("O6").Value = 0.00516242483834057
$ws.Range("P6").Value = 0.00516242483834057
$ws.Range("Q6").Value = 0.3135103077471111
$ws.Range("R6").Value = 2.821592769724
$ws.Range("S6").Value = 0.002439936090117353
$ws.Range("T6").Value = 0.002439936090117353

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Vegfc"
$ws.Range("C7").Value = "Vipr2"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 4.367310666666667
$ws.Range("H7").Value = 13.101932
$ws.Range("I7").Value = 0.4726337266929886
$ws.Range("J7").Value = 0.4726337266929886
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 6.804012333333334
$ws.Range("N7").Value = 20.412037
$ws.Range("O7").Value = 0.4893066248597758
$ws.Range("P7").Value = 0.4893066248597758
$ws.Range("Q7").Value = 29.71523563949822
$ws.Range("R7").Value = 267.437120755484
$ws.Range("S7").Value = 0.231262813603044
$ws.Range("T7").Value = 0.231262813603044

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Vegfc"
$ws.Range("C8").Value = "Vipr2"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 4.367310666666667
$ws.Range("H8").Value = 13.101932
$ws.Range("I8").Value = 0.4726337266929886
$ws.Range("J8").Value = 0.4726337266929886
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.029003
$ws.Range("N8").Value = 0.087009
$ws.Range("O8").Value = 0.002085734026566003
$ws.Range("P8").Value = 0.002085734026566003
$ws.Range("Q8").Value = 0.1266651112653333
$ws.Range("R8").Value = 1.139986001388
$ws.Range("S8").Value = 0.000985788245866263
$ws.Range("T8").Value = 0.000985788245866263

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Vegfc"
$ws.Range("C9").Value = "Vipr2"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 4.367310666666667
$ws.Range("H9").Value = 13.101932
$ws.Range("I9").Value = 0.4726337266929886
$ws.Range("J9").Value = 0.4726337266929886
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 7.000615333333333
$ws.Range("N9").Value = 21.001846
$ws.Range("O9").Value = 0.5034452162753175
$ws.Range("P9").Value = 0.5034452162753175
$ws.Range("Q9").Value = 30.57386201849689
$ws.Range("R9").Value = 275.164758166472
$ws.Range("S9").Value = 0.2379451887539609
$ws.Range("T9").Value = 0.2379451887539609

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Vegfc"
$ws.Range("C10").Value = "Vipr2"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.609937666666666
$ws.Range("H10").Value = 4.829813
$ws.Range("I10").Value = 0.1742286952351946
$ws.Range("J10").Value = 0.1742286952351946
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.07178566666666666
$ws.Range("N10").Value = 0.215357
$ws.Range("O10").Value = 0.00516242483834057
$ws.Range("P10").Value = 0.00516242483834057
$ws.Range("Q10").Value = 0.1155704486934444
$ws.Range("R10").Value = 1.040134038241
$ws.Range("S10").Value = 0.0008994425438338379
$ws.Range("T10").Value = 0.0008994425438338379

# Row 11
$ws.Range("A11").Value = "sCs"
$ws.Range("B11").Value = "Vegfc"
$ws.Range("C11").Value = "Vipr2"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 1.609937666666666
$ws.Range("H11").Value = 4.829813
$ws.Range("I11").Value = 0.1742286952351946
$ws.Range("J11").Value = 0.1742286952351946
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 6.804012333333334
$ws.Range("N11").Value = 20.412037
$ws.Range("O11").Value = 0.4893066248597758
$ws.Range("P11").Value = 0.4893066248597758
$ws.Range("Q11").Value = 10.95403573989789
$ws.Range("R11").Value = 98.586321659081
$ws.Range("S11").Value = 0.08525125481925556
$ws.Range("T11").Value = 0.08525125481925556

# Row 12
$ws.Range("A12").Value = "sCs"
$ws.Range("B12").Value = "Vegfc"
$ws.Range("C12").Value = "Vipr2"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 1.609937666666666
$ws.Range("H12").Value = 4.829813
$ws.Range("I12").Value = 0.1742286952351946
$ws.Range("J12").Value = 0.1742286952351946
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.029003
$ws.Range("N12").Value = 0.087009
$ws.Range("O12").Value = 0.002085734026566003
$ws.Range("P12").Value = 0.002085734026566003
$ws.Range("Q12").Value = 0.04669302214633333
$ws.Range("R12").Value = 0.420237199317
$ws.Range("S12").Value = 0.0003633947180562434
$ws.Range("T12").Value = 0.0003633947180562434

# Row 13
$ws.Range("A13").Value = "sCs"
$ws.Range("B13").Value = "Vegfc"
$ws.Range("C13").Value = "Vipr2"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 1.609937666666666
$ws.Range("H13").Value = 4.829813
$ws.Range("I13").Value = 0.1742286952351946
$ws.Range("J13").Value = 0.1742286952351946
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 7.000615333333333
$ws.Range("N13").Value = 21.001846
$ws.Range("O13").Value = 0.5034452162753175
$ws.Range("P13").Value = 0.5034452162753175
$ws.Range("Q13").Value = 11.27055431497755
$ws.Range("R13").Value = 101.434988834798
$ws.Range("S13").Value = 0.0877146031540489
$ws.Range("T13").Value = 0.0877146031540489
